$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
# Copy the existing header formatting (bold, centered, bordered) from B1 into A1
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Value = "Age classes"
$ws.Range("B1").Value = "% SFH"
$ws.Range("C1").Value = "% TH"
$ws.Range("D1").Value = "% AB"

# --- Row 2 (existing row, label & updated ratios) ---
$ws.Range("A2").Value = "1955 and before"
$ws.Range("B2").Value = 0.5542353624916999
$ws.Range("C2").Value = 0.1977404405643929
$ws.Range("D2").Value = 0.2480241969439072

# --- Build the remaining age-cohort rows (3-8) ---
# Copy A2's formatting (border/bold/center style) down so new A-column
# cells match the existing label-cell style.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A8").PasteSpecial(-4122) | Out-Null

$rows = @(
    @{ Row = 3; Label = "1956 - 1970";     SFH = 0.5290045118645327; TH = 0.1893335251958855; AB = 0.2816619629395818 },
    @{ Row = 4; Label = "1971 - 1980";     SFH = 0.5726899924657374; TH = 0.1890630662102659; AB = 0.2382469413239968 },
    @{ Row = 5; Label = "1981 - 1990";     SFH = 0.6122221616745076; TH = 0.2137909529604805; AB = 0.1739868853650119 },
    @{ Row = 6; Label = "1991 - 2000";     SFH = 0.476144176678936;  TH = 0.2414751024376693; AB = 0.2823807208833947 },
    @{ Row = 7; Label = "2001 - 2010";     SFH = 0.336760881642224;  TH = 0.2304784821084705; AB = 0.4327606362493056 },
    @{ Row = 8; Label = "2011 and after";  SFH = 0.340234951214744;  TH = 0.2475518509961434; AB = 0.4122131977891126 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Label
    $ws.Cells.Item($r.Row, 2).Value = $r.SFH
    $ws.Cells.Item($r.Row, 3).Value = $r.TH
    $ws.Cells.Item($r.Row, 4).Value = $r.AB
}

Write-Host "done"
